# "Generate Report for Archive"
#
# The localization status report is being refreshed: items that were
# previously marked "Ready for handoff" have moved on to "In Translation".
# That status string shows up in three places:
#   - Overview sheet, zh-cn (col E) and de-de (col F) status columns, rows 2-4
#   - zh-cn sheet, Status column (col C), rows 2-4
#   - de-de sheet, Status column (col C), rows 2-4
#
# Because "In Translation" is shorter than "Ready for handoff", the status
# columns are narrowed to fit the new text.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

$oldStatus = "Ready for handoff"
$newStatus = "In Translation"

# --- Update the status cells ---
$wsOverview.Range("E2:F4").Value = $newStatus
$wsZhCn.Range("C2:C4").Value = $newStatus
$wsDeDe.Range("C2:C4").Value = $newStatus

# --- Re-fit the status columns now that the text is shorter ---
$wsOverview.Range("E1:F1").ColumnWidth = 12.5
$wsZhCn.Range("C1").ColumnWidth = 12.5
$wsDeDe.Range("C1").ColumnWidth = 12.5
